# Update attendance summary cells (Total Attendance Count / Real / Invalid / Absent
# columns) from 0 to 1 for the rows/columns that reflect a recorded attendance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToSet = @(
    "G3", "H3",
    "D4", "E4",
    "D5", "E5",
    "D6", "E6",
    "H7",
    "H8",
    "D9", "E9",
    "G10", "H10",
    "D11", "E11",
    "D12", "E12",
    "H13",
    "H14",
    "D15", "E15",
    "H16",
    "D17", "E17",
    "H18"
)

foreach ($addr in $cellsToSet) {
    $ws.Range($addr).Value = 1
}
